$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B7").Value = 3935283.45
$ws.Range("C7").Value = -11.42889357975097
$ws.Range("D7").Value = 3406
$ws.Range("E7").Value = 3406
$ws.Range("F7").Value = 1155.397372284204
$ws.Range("G7").Value = 23.15700528664106
